# Clean up stale/duplicate quantity figures and fix inline-string numeric
# cells on the stock data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Decrement a handful of quantity values (column D) that had drifted out
# of sync with the real stock counts.
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 43
$ws.Range("D7").Value = 44
$ws.Range("D9").Value = 6
$ws.Range("D11").Value = 9

# Row 14 ("Papel Contact Pliego") had its numbers stored as text; rewrite
# them as real numbers, with the quantity corrected at the same time.
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = 400
$ws.Range("F14").Value = 650
